$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1 ("Metadata"): insert a new "Jurisdiction" row right after "Contact"
# (row 10), pushing Description/Purpose/Copyright/Immutable down by one row,
# and update the Version and Date values.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Shift rows 11-14 down to 12-15 (process bottom-up so we don't clobber data).
# Copying the source row's formatting (xlPasteFormats) onto the destination
# row before writing the values keeps the existing cell style (border/
# alignment) instead of acquiring a brand new style index.
for ($r = 14; $r -ge 11; $r--) {
    $colA = $ws1.Cells.Item($r, 1).Value2
    $colB = $ws1.Cells.Item($r, 2).Value2

    $ws1.Range("A" + $r + ":B" + $r).Copy()
    $ws1.Range("A" + ($r+1) + ":B" + ($r+1)).PasteSpecial(-4122)

    $ws1.Cells.Item($r+1, 1).Value = $colA
    if ($colB -eq $null) {
        $ws1.Cells.Item($r+1, 2).Value = ""
    } else {
        $ws1.Cells.Item($r+1, 2).Value = $colB
    }
}
$excel.CutCopyMode = $false

# New row 11: Jurisdiction / iso:code:3166:FR
$ws1.Cells.Item(11, 1).Value = "Jurisdiction"
$ws1.Cells.Item(11, 2).Value = "iso:code:3166:FR"

# Update Version (row 3) and Date (row 8)
$ws1.Cells.Item(3, 2).Value = "0.2.0"
$ws1.Cells.Item(8, 2).Value = "2023-10-20T08:59:58+00:00"

# ---------------------------------------------------------------------------
# Sheet 2 ("Include from Type de contact"): unchanged content, indices in
# the shared string table simply move because two new strings were added
# above - nothing to do here, values stay the same.
# ---------------------------------------------------------------------------
